$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Infused Beverages"
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Suspicious Brew"
Write-Output $wb.Worksheets.Count
Write-Output $wb.Worksheets.Item(1).Name
Write-Output $wb.Worksheets.Item(2).Name
